$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.445.35'
$ws.Range('E2').Value = '  -1.03%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.906.32'
$ws.Range('E3').Value = '  +3.64%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.35'
$ws.Range('E5').Value = '  +0.27%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.79'
$ws.Range('E6').Value = '  -0.91%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.904.67'
$ws.Range('E7').Value = '  +3.65%  '

$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('E9').Value = '  -2.35%  '

$ws.Range('E10').Value = '  -3.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.39'
$ws.Range('E11').Value = '  +0.24%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.83'
$ws.Range('E13').Value = '  -2.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000245'
$ws.Range('E14').Value = '  -1.08%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.556.71'
$ws.Range('E15').Value = '  +3.68%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.890.47'
$ws.Range('E16').Value = '  +3.18%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.652.41'
$ws.Range('E17').Value = '  -0.87%  '

$ws.Range('E18').Value = '  -0.07%  '

$ws.Range('E19').Value = '  -1.20%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.01'
$ws.Range('E20').Value = '  -3.80%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.16'
$ws.Range('E21').Value = '  -0.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '483.64'
$ws.Range('E22').Value = '  -1.78%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.717'
$ws.Range('E23').Value = '  -1.21%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000168'
$ws.Range('E24').Value = '  +12.71%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.28'
$ws.Range('E25').Value = '  -0.60%  '

$ws.Range('E26').Value = '  -1.58%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.00'
$ws.Range('E27').Value = '  -2.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  +0.23%  '

$ws.Range('E29').Value = '  +0.01%  '

$ws.Range('E30').Value = '  -1.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.054.12'
$ws.Range('E31').Value = '  +3.74%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.85'
$ws.Range('E32').Value = '  -3.27%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.37'
$ws.Range('E33').Value = '  -2.40%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.97'
$ws.Range('E34').Value = '  +0.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.847.26'
$ws.Range('E35').Value = '  +3.35%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.107'
$ws.Range('E36').Value = '  -0.91%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  +2.44%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.140'
$ws.Range('E38').Value = '  +0.90%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.87'
$ws.Range('E39').Value = '  -1.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.05'
$ws.Range('E41').Value = '  -1.75%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.317'
$ws.Range('E42').Value = '  -2.14%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '433.55'
$ws.Range('E43').Value = '  +1.87%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.47'
$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.98'
$ws.Range('E45').Value = '  -0.77%  '

$ws.Range('E46').Value = '  +0.01%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.42'
$ws.Range('E47').Value = '  -0.30%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.35'
$ws.Range('E48').Value = '  +11.60%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.814.64'
$ws.Range('E49').Value = '  -0.14%  '

$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.56'
$ws.Range('E50').Value = '  -0.51%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.40'
$ws.Range('E51').Value = '  -1.91%  '

